# Update "想去人数" (wanted-to-go count) figures for the 2024-05-01 exhibition
# and the 2024-06-22 performance across the relevant sheets.

$wb = $excel.ActiveWorkbook

# "展览" (Exhibition) sheet - F2 holds the count for the 2024-05-01 event
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 5476

# "演出" (Performance) sheet - F2 holds the count for the 2024-06-22 event
$wsPerformance = $wb.Worksheets.Item("演出")
$wsPerformance.Range("F2").Value = 3

# "全部类型" (All types) sheet - aggregates both events above
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 5476
$wsAll.Range("F5").Value = 3
